$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) for new columns I and J - copy formatting from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for rows 2-70, columns I (I0) and J (IF)
$data = @(
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(6,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(9,9),
    @(7,8),
    @(8,8),
    @(8,8),
    @(7,8),
    @(6,7),
    @(7,8),
    @(7,8),
    @(8,8),
    @(5,6),
    @(7,7),
    @(9,9),
    @(5,5),
    @(9,9),
    @(6,7),
    @(3,4),
    @(5,6),
    @(6,6),
    @(9,9),
    @(9,9),
    @(8,9),
    @(3,4),
    @(10,10),
    @(9,9),
    @(5,5),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,9),
    @(11,11),
    @(7,8),
    @(5,6),
    @(6,6),
    @(7,8),
    @(7,7),
    @(9,9),
    @(6,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(6,7),
    @(7,7),
    @(8,9),
    @(8,8),
    @(8,8),
    @(6,7),
    @(8,8),
    @(7,8),
    @(9,9),
    @(8,9),
    @(7,9),
    @(5,5),
    @(6,6),
    @(1,2),
    @(7,7),
    @(5,5),
    @(6,7),
    @(7,8)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $data[$r][0]
    $ws.Cells.Item($row, 10).Value = $data[$r][1]
}
